# The data rows 2-10 (one row per atomA/atomB/atomC combination) need to have
# their "identity + output coordinate" columns (A, B, C and S through AD)
# reversed in row order -- i.e. row 2 swaps with row 10, row 3 swaps with
# row 9, row 4 swaps with row 8, row 5 swaps with row 7, and row 6 (the
# middle row) stays put. The other columns (D through R) are left untouched
# because they stay associated with their original row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together with the swapped rows.
$cols = @("A", "B", "C", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD")

# Row pairs (1-based worksheet rows) that need to trade their values.
$pairs = @(
    @(2, 10),
    @(3, 9),
    @(4, 8),
    @(5, 7)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        # NOTE: `.Value` doesn't invoke the COM getter correctly in this
        # host - use `.Value2` for both reading and writing instead.
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
